# octanol 30pct avoidance - append 09/30 data (row 57 onward)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$odor          = "30pct"
$experimenter  = "TB"
$expDate       = 44468   # 2021-09-30 (serial date number)

# ---- N2, with odor recorded (rows 57-78) ----
$n2Times = @(2.76,7.93,6.47,3.71,3.23,5.73,3.94,4.47,6.68,9.25,5.6, `
             6.37,2.15,7.39,6.13,6.04,4.83,7.2,4.22,5.25,3.53,5.55)

$row = 57
foreach ($t in $n2Times) {
  $ws.Cells.Item($row,1).Value = "N2"
  $ws.Cells.Item($row,2).Value = $odor
  $ws.Cells.Item($row,3).Value = $t
  $ws.Cells.Item($row,4).Value = $experimenter
  $ws.Cells.Item($row,5).Value = $expDate
  $row++
}

# ---- cest-2.1, first entry still has odor recorded (row 79) ----
$ws.Cells.Item($row,1).Value = "cest-2.1"
$ws.Cells.Item($row,2).Value = $odor
$ws.Cells.Item($row,3).Value = 1.15
$ws.Cells.Item($row,4).Value = $experimenter
$ws.Cells.Item($row,5).Value = $expDate
$row++

# ---- cest-2.1, remaining entries - odor column left blank (rows 80-101) ----
$cestTimes = @(6.56,5.16,3.25,4.23,3.84,5.68,3.21,5.11,2.06,1.51,2.71, `
               3.96,1.62,2.84,1.85,2.96,8.38,3.88,2.2,3.57,4.09,6.27)

foreach ($t in $cestTimes) {
  $ws.Cells.Item($row,1).Value = "cest-2.1"
  $ws.Cells.Item($row,3).Value = $t
  $ws.Cells.Item($row,4).Value = $experimenter
  $ws.Cells.Item($row,5).Value = $expDate
  $row++
}

$lastRow = $row - 1

# ---- formatting ----
# Give the new date cells (E57:E<lastRow>) the same style the rest of
# column E already uses (date format), reusing the existing style so we
# don't fragment the stylesheet.
$ws.Range("E2").Copy() | Out-Null
$ws.Range("E57:E" + $lastRow).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Re-apply the date number format across the whole date column (old +
# new rows) as yyyy-mm-dd instead of the old m/d/yyyy.
$ws.Range("E2:E" + $lastRow).NumberFormat = "yyyy\-mm\-dd;@"

# Column E now holds data, give it a sensible width like the other
# data columns (C, D) already have.
$ws.Columns.Item(5).ColumnWidth = 10

# Printed page: portrait like the rest of the workbook.
$ws.PageSetup.Orientation = 1

# ---- view state ----
# Leave the selection where the user's data entry ended, and scroll the
# window down so the new rows are visible.
$excel.ActiveWindow.ScrollRow = 92
$ws.Range("D105").Select() | Out-Null
